# Atualização automática de CERRO_LARGO.xlsx
#
# 1) Rename "Paineis DARQ" -> "PAINEIS DARQ"
# 2) Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3) Delete the "Desarquivamentos Pendentes" sheet (the CERRO LARGO
#    pending-requests panel) entirely.
# 4) Leave the first sheet ("PAINEIS DARQ") selected/active, matching
#    the original workbook's tab selection.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
[void]$wb.Worksheets("Desarquivamentos Pendentes").Delete()

[void]$wb.Worksheets("PAINEIS DARQ").Activate()

$excel.DisplayAlerts = $true
